$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "[[-0.01259479]`n [ 0.03034212]`n [ 0.06248812]`n [ 0.02095435]`n [-0.05816003]`n [-0.1876165 ]`n [-0.0685947 ]`n [ 0.03341826]`n [ 0.14568317]`n [ 0.0732327 ]`n [ 0.0585992 ]`n [ 0.08114034]`n [-0.22276638]]"
$ws.Range("E2").Value = "[[-0.00512007]`n [-0.02061025]`n [ 0.15090641]`n [ 0.01116715]`n [ 0.07244902]`n [-0.16741335]`n [ 0.10327631]`n [ 0.05800829]`n [ 0.08898542]`n [ 0.05498782]`n [-0.05554802]`n [-0.07891823]`n [-0.25604462]]"
$ws.Range("F2").Value = "[-0.00728113  0.03071338 -0.07319244 -0.01381209 -0.13809748 -0.02902286`n -0.09674671  0.01028812  0.03443562  0.04515991  0.09864035 -0.16041014`n -0.02462132]"
$ws.Range("G2").Value = "[-0.00929498  0.03179126  0.05999124  0.02234818 -0.05756453 -0.18848256`n -0.07088835  0.03299387  0.14752245  0.0740176   0.06037279  0.08233421`n -0.22901532]"
$ws.Range("D3").Value = "[[ 0.1021785 ]`n [-0.01395554]`n [-0.04214394]`n [-0.10800695]`n [ 0.10941466]`n [ 0.01154953]`n [ 0.06453123]`n [-0.01845728]`n [-0.03349375]`n [ 0.00092294]`n [-0.01551101]`n [-0.10548255]`n [ 0.00458001]]"
$ws.Range("E3").Value = "[[ 0.06842089]`n [-0.0793725 ]`n [-0.05879938]`n [-0.077634  ]`n [ 0.08275873]`n [-0.03398016]`n [ 0.08914596]`n [-0.04155097]`n [ 0.0043805 ]`n [ 0.0708936 ]`n [-0.00629081]`n [-0.07195538]`n [ 0.01010939]]"
$ws.Range("F3").Value = "[ 0.05118617 -0.08705603 -0.01671972 -0.05980253  0.02321769  0.03433534`n  0.01032242  0.0184057  -0.12173654 -0.05782419  0.00099413 -0.14246176`n -0.01087969]"
$ws.Range("G3").Value = "[ 0.09550865 -0.01066191 -0.03435089 -0.10577454  0.10775223  0.00708206`n  0.06304329 -0.01701662 -0.03901747  0.00908708 -0.01760989 -0.1065105`n  0.00459438]"
